$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# inner_sheet_thick (A2): 0.6 -> 0.8
$ws.Range("A2").Value = 0.8

# t1_length (C2): 400 -> 350
$ws.Range("C2").Value = 350

# t1_height (D2): 300 -> 200
$ws.Range("D2").Value = 200

# panel_thick (K2): 50 -> 48
$ws.Range("K2").Value = 48
